# Update column G ("K") values on the active sheet to reflect the
# regenerated Strike# -> K calculation (std/mean based s_vals).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 6
    3  = 2
    4  = 6
    5  = 6
    6  = 8
    7  = 5
    8  = 7
    9  = 7
    10 = 6
    11 = 7
    12 = 6
    13 = 5
    14 = 8
    15 = 7
    16 = 5
    17 = 4
    18 = 6
    19 = 3
    20 = 5
    21 = 5
    22 = 7
    23 = 5
    24 = 2
    25 = 5
    26 = 2
    27 = 2
    28 = 4
    29 = 8
    30 = 4
    31 = 4
    32 = 5
    33 = 10
    34 = 7
    35 = 5
    36 = 8
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
